$wb = $excel.ActiveWorkbook
$county = $wb.Worksheets.Item("County")
$county.Range("B1").Value = 170
